# Rename the first sheet from "Plan1" to "Cripto"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Cripto"

# Write the header row used to track cryptocurrency quotes
$ws.Range("A1").Value = "Criptomoeda"
$ws.Range("B1").Value = "Valor da última consulta"
$ws.Range("C1").Value = "Valor do dia atual"
$ws.Range("D1").Value = "% de aumento"
$ws.Range("E1").Value = "% de queda"

# Bold the header
$ws.Range("A1:E1").Font.Bold = $true

# Size the columns so the header text fits (best-fit)
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 21.833333333333332
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 10.5

# Touch the page setup so printing defaults (A4/portrait) are persisted
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor where the user stopped entering data
$ws.Range("F14").Select() | Out-Null
